# Update Name of Algo
# Apply updated values to column A and D per the new algorithm run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -22.17500000000001
$ws.Range("A21").Value = -20.02119999999998
$ws.Range("A23").Value = -20.29019999999998
$ws.Range("D24").Value = -7.6189
$ws.Range("A25").Value = -21.75
$ws.Range("D28").Value = -7.999199999999997
$ws.Range("D36").Value = -7.0905
$ws.Range("D45").Value = -7.101900000000001
$ws.Range("D48").Value = -7.447799999999993
$ws.Range("D49").Value = -8.041700000000004
$ws.Range("D52").Value = -7.837600000000004
$ws.Range("A53").Value = -22.1156
$ws.Range("D53").Value = -8.171799999999998
$ws.Range("D54").Value = -8.105600000000003
$ws.Range("A57").Value = -22.58570000000001
$ws.Range("A59").Value = -22.14489999999999
$ws.Range("A69").Value = -21.56359999999999
$ws.Range("D70").Value = -6.622499999999999
$ws.Range("A79").Value = -19.88489999999999
$ws.Range("A83").Value = -21.79340000000001
$ws.Range("D86").Value = -8.748400000000006
$ws.Range("D87").Value = -8.422499999999992
$ws.Range("A93").Value = -21.51320000000002
$ws.Range("D101").Value = -7.910700000000002
